$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking strings
# (e.g. "1.00", "0.999") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.143.27"
$ws.Range("E2").Value = "  +5.44%  "
$ws.Range("D3").Value = "2.259.14"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "301.90"
$ws.Range("E5").Value = "  +3.41%  "
$ws.Range("D6").Value = "93.11"
$ws.Range("E6").Value = "  +6.23%  "
$ws.Range("D7").Value = "0.533"
$ws.Range("E7").Value = "  +3.78%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.484"
$ws.Range("E9").Value = "  +3.59%  "
$ws.Range("D10").Value = "32.80"
$ws.Range("E10").Value = "  +7.76%  "
$ws.Range("D11").Value = "54.59"
$ws.Range("E11").Value = "  +8.93%  "
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").Value = "6.68"
$ws.Range("E14").Value = "  +3.79%  "
$ws.Range("D15").Value = "2.607.07"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").Value = "2.249.33"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").Value = "0.758"
$ws.Range("E18").Value = "  +3.69%  "
$ws.Range("D19").Value = "41.989.05"
$ws.Range("E19").Value = "  +5.20%  "
$ws.Range("D20").Value = "12.24"
$ws.Range("E20").Value = "  +10.03%  "
$ws.Range("E21").Value = "  +2.08%  "
$ws.Range("D22").Value = "5.96"
$ws.Range("E22").Value = "  +3.85%  "
$ws.Range("D23").Value = "67.30"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("D24").Value = "242.35"
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +5.27%  "
$ws.Range("E27").Value = "  +4.61%  "
$ws.Range("D28").Value = "24.05"
$ws.Range("E28").Value = "  +3.81%  "
$ws.Range("D29").Value = "9.72"
$ws.Range("E29").Value = "  +5.07%  "
$ws.Range("D30").Value = "2.08"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "159.33"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "34.22"
$ws.Range("E32").Value = "  +7.04%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "5.17"
$ws.Range("E34").Value = "  +4.08%  "
$ws.Range("D35").Value = "0.0746"
$ws.Range("E35").Value = "  +4.58%  "
$ws.Range("E36").Value = "  +3.13%  "
$ws.Range("E37").Value = "  +2.48%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "16.86"
$ws.Range("E38").Value = "  +10.14%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.105"
$ws.Range("E39").Value = "  +5.80%  "
$ws.Range("E40").Value = "  +3.77%  "
$ws.Range("D41").Value = "1.81"
$ws.Range("E41").Value = "  +4.48%  "
$ws.Range("D42").Value = "3.95"
$ws.Range("E42").Value = "  +5.83%  "
$ws.Range("D43").Value = "2.055.35"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("D44").Value = "20.03"
$ws.Range("E44").Value = "  +11.98%  "
$ws.Range("E45").Value = "  +3.76%  "
$ws.Range("D46").Value = "10.13"
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("D47").Value = "2.89"
$ws.Range("E47").Value = "  +7.05%  "
$ws.Range("D48").Value = "2.04"
$ws.Range("E48").Value = "  +4.03%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "1.53"
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").Value = "  +4.26%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "51.96"
$ws.Range("E51").Value = "  +5.91%  "
